# Edit script: insert 4 new data rows (a new reporting date, 2021-10-22 / serial 44491)
# into the "Choclo" sheet right before the current row 373, pushing the existing
# rows 373:431 down to 377:435. This matches the commit "Fruta / hortaliza, semanal"
# which adds a new week of price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 373 (existing rows 373-431 shift down to 377-435)
$ws.Range("A373:R376").EntireRow.Insert()

# New row 373: Choclo, Dulce o Americano, Primera
$ws.Cells.Item(373, 1).Value = 1
$ws.Cells.Item(373, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(373, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(373, 4).Value = 44491
$ws.Cells.Item(373, 5).Value = 15
$ws.Cells.Item(373, 6).Value = 100112024
$ws.Cells.Item(373, 7).Value = "Choclo"
$ws.Cells.Item(373, 8).Value = "Dulce o Americano"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 100
$ws.Cells.Item(373, 11).Value = 26000
$ws.Cells.Item(373, 12).Value = 27000
$ws.Cells.Item(373, 13).Value = 26500
$ws.Cells.Item(373, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(373, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(373, 16).Value = 379
$ws.Cells.Item(373, 17).Value = 70
$ws.Cells.Item(373, 18).Value = "Hortaliza"

# New row 374: Choclo, Lluteño, Primera
$ws.Cells.Item(374, 1).Value = 1
$ws.Cells.Item(374, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(374, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(374, 4).Value = 44491
$ws.Cells.Item(374, 5).Value = 15
$ws.Cells.Item(374, 6).Value = 100112024
$ws.Cells.Item(374, 7).Value = "Choclo"
$ws.Cells.Item(374, 8).Value = "Lluteño"
$ws.Cells.Item(374, 9).Value = "Primera"
$ws.Cells.Item(374, 10).Value = 40
$ws.Cells.Item(374, 11).Value = 31000
$ws.Cells.Item(374, 12).Value = 32000
$ws.Cells.Item(374, 13).Value = 31500
$ws.Cells.Item(374, 14).Value = "`$/saco 50 unidades"
$ws.Cells.Item(374, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(374, 16).Value = 630
$ws.Cells.Item(374, 17).Value = 50
$ws.Cells.Item(374, 18).Value = "Hortaliza"

# New row 375: Choclo, Lluteño, Segunda
$ws.Cells.Item(375, 1).Value = 1
$ws.Cells.Item(375, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(375, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(375, 4).Value = 44491
$ws.Cells.Item(375, 5).Value = 15
$ws.Cells.Item(375, 6).Value = 100112024
$ws.Cells.Item(375, 7).Value = "Choclo"
$ws.Cells.Item(375, 8).Value = "Lluteño"
$ws.Cells.Item(375, 9).Value = "Segunda"
$ws.Cells.Item(375, 10).Value = 50
$ws.Cells.Item(375, 11).Value = 29000
$ws.Cells.Item(375, 12).Value = 30000
$ws.Cells.Item(375, 13).Value = 29500
$ws.Cells.Item(375, 14).Value = "`$/saco 75 unidades"
$ws.Cells.Item(375, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(375, 16).Value = 393
$ws.Cells.Item(375, 17).Value = 75
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# New row 376: Choclo, Lluteño, Tercera
$ws.Cells.Item(376, 1).Value = 1
$ws.Cells.Item(376, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(376, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(376, 4).Value = 44491
$ws.Cells.Item(376, 5).Value = 15
$ws.Cells.Item(376, 6).Value = 100112024
$ws.Cells.Item(376, 7).Value = "Choclo"
$ws.Cells.Item(376, 8).Value = "Lluteño"
$ws.Cells.Item(376, 9).Value = "Tercera"
$ws.Cells.Item(376, 10).Value = 60
$ws.Cells.Item(376, 11).Value = 24000
$ws.Cells.Item(376, 12).Value = 25000
$ws.Cells.Item(376, 13).Value = 24500
$ws.Cells.Item(376, 14).Value = "`$/saco 100 unidades"
$ws.Cells.Item(376, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(376, 16).Value = 245
$ws.Cells.Item(376, 17).Value = 100
$ws.Cells.Item(376, 18).Value = "Hortaliza"

# Ensure the date cells use the same date/time number format as the rest of column D
$ws.Range("D373:D376").NumberFormat = $ws.Range("D372").NumberFormat
